$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> Jason Roy)
$ws.Name = "Jason Roy"

# Header row
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")

# Data rows (matchNo, teamName, batterName, states, runs, balls, fours, sixes, sr, opponentTeamName, venue, date, result)
$rows = @(
    @("44th","Sunrisers Hyderabad","Jason Roy","c †Dhoni b Hazlewood","2","7","0","0","28.57","Chennai Super Kings","Sharjah","September 30","Super Kings won by 6 wickets (with 2 balls remaining)"),
    @("49th","Sunrisers Hyderabad","Jason Roy","c Southee b Shivam Mavi","10","13","2","0","76.92","Kolkata Knight Riders","Dubai (DSC)","October 03","KKR won by 6 wickets (with 2 balls remaining)"),
    @("52nd","Sunrisers Hyderabad","Jason Roy","c & b Christian","44","38","5","0","115.78","Royal Challengers Bangalore","Abu Dhabi","October 06","Sunrisers won by 4 runs"),
    @("55th","Sunrisers Hyderabad","Jason Roy","c KH Pandya b Boult","34","21","6","0","161.90","Mumbai Indians","Abu Dhabi","October 08","Mumbai won by 42 runs"),
    @("40th","Sunrisers Hyderabad","Jason Roy","c †Samson b Sakariya","60","42","8","1","142.85","Rajasthan Royals","Dubai (DSC)","September 27","Sunrisers won by 7 wickets (with 9 balls remaining)")
)

$lastCol = $headers.Length
$lastRow = 1 + $rows.Length

# Force the whole used range to Text so numeric-looking strings (e.g. "2", "28.57")
# are stored as text, matching the source data's string typing.
$usedRange = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item($lastRow, $lastCol))
$usedRange.NumberFormat = "@"

for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}
